# DGS risk-matrix time series update — add the 2021/08/20 report row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds dates as text (formatted yyyy/mm/dd); a leading apostrophe
# forces the new entry to be stored as text too, matching the existing rows,
# instead of letting Excel auto-convert it to a date serial number.
$ws.Range("A70").Value = "'2021/08/20"
$ws.Range("B70").Value = 312.3
$ws.Range("C70").Value = 316.6
$ws.Range("D70").Value = 0.98
$ws.Range("E70").Value = 0.98

# Move the selection to the next empty row, as happens after data entry.
[void]$ws.Range("A71").Select()
